$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 and 3 hold two weekly price records that need to swap places
# (row 2 should end up with the later date's data, row 3 with the
# earlier date's data) for columns D, J, K, L, M, P.

$D2 = $ws.Range("D2").Value2
$J2 = $ws.Range("J2").Value2
$K2 = $ws.Range("K2").Value2
$L2 = $ws.Range("L2").Value2
$M2 = $ws.Range("M2").Value2
$P2 = $ws.Range("P2").Value2

$D3 = $ws.Range("D3").Value2
$J3 = $ws.Range("J3").Value2
$K3 = $ws.Range("K3").Value2
$L3 = $ws.Range("L3").Value2
$M3 = $ws.Range("M3").Value2
$P3 = $ws.Range("P3").Value2

$ws.Range("D2").Value2 = $D3
$ws.Range("J2").Value2 = $J3
$ws.Range("K2").Value2 = $K3
$ws.Range("L2").Value2 = $L3
$ws.Range("M2").Value2 = $M3
$ws.Range("P2").Value2 = $P3

$ws.Range("D3").Value2 = $D2
$ws.Range("J3").Value2 = $J2
$ws.Range("K3").Value2 = $K2
$ws.Range("L3").Value2 = $L2
$ws.Range("M3").Value2 = $M2
$ws.Range("P3").Value2 = $P2
